# Event Hawk presentation - "Changes from our presentation"
#
# Content edits:
#   1. Slide 3 ("Current status"): bullet "No backend integration"
#      -> "Initial backend integration"
#   2. Slide 4 ("Goals"): bullet "Integrate with backend"
#      -> "Complete integration with backend"
#   3. The cached "datetimeFigureOut" date field text (10/10/2017) on the
#      slide master and every slide layout is refreshed to 10/17/2017.

$p = $ppt.ActivePresentation

# --- 1. Slide 3 : "Current status" ---------------------------------------
$s3 = $p.Slides.Item(3)
$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange
for ($i = 1; $i -le $s3Body.Paragraphs().Count; $i++) {
    $para = $s3Body.Paragraphs($i)
    $run = $para.Runs(1)
    if ($run.Text -eq "No backend integration") {
        $run.Text = "Initial backend integration"
    }
}

# --- 2. Slide 4 : "Goals" --------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
for ($i = 1; $i -le $s4Body.Paragraphs().Count; $i++) {
    $para = $s4Body.Paragraphs($i)
    $run = $para.Runs(1)
    if ($run.Text -eq "Integrate with backend") {
        $run.Text = "Complete integration with backend"
    }
}

# --- 3. Refresh the cached date field (10/10/2017 -> 10/17/2017) ----------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "10/10/2017") {
                $tr.Text = "10/17/2017"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
